$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New roster table (header row 1 stays as-is)
$data = @(
    @("Immanuel Quickley", "PG,SG",      "Toronto Raptors"),
    @("Norman Powell",     "SG,SF",      "LA Clippers"),
    @("Christian Braun",   "SG,SF",      "Denver Nuggets"),
    @("Jimmy Butler",      "SF,PF",      "Golden State Warriors"),
    @("Kawhi Leonard",     "SG,SF,PF",   "LA Clippers"),
    @("Jalen Williams",    "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Jalen Brunson",     "PG",         "New York Knicks"),
    @("Corey Kispert",     "SG,SF",      "Washington Wizards"),
    @("Desmond Bane",      "SG,SF",      "Memphis Grizzlies"),
    @("Devin Booker",      "PG,SG",      "Phoenix Suns"),
    @("Trae Young",        "PG",         "Atlanta Hawks"),
    @("Walker Kessler",    "C",          "Utah Jazz"),
    @("LeBron James",      "SF,PF",      "Los Angeles Lakers"),
    @("Tyus Jones",        "PG",         "Phoenix Suns"),
    @("Brandon Ingram",    "SG,SF,PF",   "Toronto Raptors"),
    @("Myles Turner",      "C",          "Indiana Pacers")
)

# Remove the two rows no longer present (table shrinks from 18 to 16 data rows)
$ws.Range("A18:C19").Delete()

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
